$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 12; $r++) {
    $ws.Range("AB${r}:AK${r}").ClearContents()
    $ws.Range("AM${r}").ClearContents()
}
